$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; everything (values + formatting) shifts down by one row.
$ws.Rows.Item(1).Insert()

# The old header row ("Lg.", "Threading", ... with style s=1) is now row 2.
# Copy its formatting up onto the brand-new row 1 (reuses the existing style, e.g. s="1"),
# then strip the inherited formatting off row 2 so it goes back to the default style.
$ws.Range("A2:L2").Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)
$ws.Range("A2:L2").ClearFormats()

# New row 1 becomes a simple numeric column-index header: 0,1,2,...,11
for ($i = 0; $i -lt 12; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $i
}

# Row 2 keeps the original textual headers in columns A,B,C,D,E,F,G,H,J,
# but I2/K2/L2 are cleared out (no longer carry a value).
$ws.Range("I2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
